$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Which films did Stanley Kubrick direct?"
$ws.Range("B2").Value = "['http://dbpedia.org/resource/2001:_A_Space_Odyssey_(film)', 'http://dbpedia.org/resource/A_Clockwork_Orange_(film)', 'http://dbpedia.org/resource/Barry_Lyndon', 'http://dbpedia.org/resource/Day_of_the_Fight', 'http://dbpedia.org/resource/Dr._Strangelove', 'http://dbpedia.org/resource/Eyes_Wide_Shut', 'http://dbpedia.org/resource/Fear_and_Desire', 'http://dbpedia.org/resource/Flying_Padre', 'http://dbpedia.org/resource/Full_Metal_Jacket', `"http://dbpedia.org/resource/Killer's_Kiss`", 'http://dbpedia.org/resource/Lolita_(1962_film)', 'http://dbpedia.org/resource/Paths_of_Glory', 'http://dbpedia.org/resource/Spartacus_(film)', 'http://dbpedia.org/resource/The_Killing_(film)__The_Killing__1', 'http://dbpedia.org/resource/The_Seafarers__The_Seafarers__1', 'http://dbpedia.org/resource/The_Shining_(film)']"
$ws.Range("C2").Value = "['http://dbpedia.org/resource/2001:_A_Space_Odyssey_(film)', 'http://dbpedia.org/resource/A_Clockwork_Orange_(film)', 'http://dbpedia.org/resource/Barry_Lyndon', 'http://dbpedia.org/resource/Day_of_the_Fight', 'http://dbpedia.org/resource/Dr._Strangelove', 'http://dbpedia.org/resource/Eyes_Wide_Shut', 'http://dbpedia.org/resource/Fear_and_Desire', 'http://dbpedia.org/resource/Flying_Padre', 'http://dbpedia.org/resource/Full_Metal_Jacket', `"http://dbpedia.org/resource/Killer's_Kiss`", 'http://dbpedia.org/resource/Lolita_(1962_film)', 'http://dbpedia.org/resource/Paths_of_Glory', 'http://dbpedia.org/resource/Spartacus_(film)', 'http://dbpedia.org/resource/The_Killing_(film)__The_Killing__1', 'http://dbpedia.org/resource/The_Seafarers__The_Seafarers__1', 'http://dbpedia.org/resource/The_Shining_(film)']"
$ws.Range("D2").Value = $true

$ws.Range("A3").Value = "How many films did Leonardo DiCaprio star in?"
$ws.Range("B3").Value = "['33']"
$ws.Range("C3").Value = "['33']"
$ws.Range("D3").Value = $true

$ws.Range("A4").Value = "How high is the Yokohama Marine Tower?"
$ws.Range("B4").Value = "['106.07']"
$ws.Range("C4").Value = "['106.07']"
$ws.Range("D4").Value = $true

$ws.Range("A5").Value = "What is Donald Trump's main business?"
$ws.Range("B5").Value = "['http://dbpedia.org/resource/GoTrump.com', 'http://dbpedia.org/resource/The_Trump_Organization', 'http://dbpedia.org/resource/Trump_Model_Management']"
$ws.Range("C5").Value = "['http://dbpedia.org/resource/GoTrump.com', 'http://dbpedia.org/resource/The_Trump_Organization', 'http://dbpedia.org/resource/Trump_Model_Management']"
$ws.Range("D5").Value = $true

$ws.Range("A6").Value = "Who is the author of the interpretation of dreams?"
$ws.Range("B6").Value = "['http://dbpedia.org/resource/Sigmund_Freud']"
$ws.Range("C6").Value = "['http://dbpedia.org/resource/Sigmund_Freud']"
$ws.Range("D6").Value = $true

$ws.Range("A7").Value = "What country is Mount Everest in?"
$ws.Range("B7").Value = "['http://dbpedia.org/resource/China', 'http://dbpedia.org/resource/Nepal']"
$ws.Range("C7").Value = "['http://dbpedia.org/resource/China', 'http://dbpedia.org/resource/Nepal']"
$ws.Range("D7").Value = $true

$ws.Range("A8").Value = "Which people were born in Heraklion?"
$ws.Range("B8").Value = "['http://dbpedia.org/resource/Anastasios_Triantafyllou_(footballer)', 'http://dbpedia.org/resource/Anastasios_Triantafyllou', 'http://dbpedia.org/resource/Andreas_Bouchalakis', 'http://dbpedia.org/resource/Andreas_Musalus', 'http://dbpedia.org/resource/Antonis_Alexakis', 'http://dbpedia.org/resource/Antonis_Bourselis', 'http://dbpedia.org/resource/Archbishop_Makarios_of_Australia', 'http://dbpedia.org/resource/Beata_Kitsikis', 'http://dbpedia.org/resource/Charalampos_Charalampakis', 'http://dbpedia.org/resource/Charalampos_Damianakis', 'http://dbpedia.org/resource/Christoforos_Liontakis', 'http://dbpedia.org/resource/Christoforos_Papakaliatis', 'http://dbpedia.org/resource/Christos_Belevonis', 'http://dbpedia.org/resource/Christos_Kozyrakis', 'http://dbpedia.org/resource/Christos_Saloustros', 'http://dbpedia.org/resource/Chrysovalantis_Kozoronis', 'http://dbpedia.org/resource/Cyril_Lucaris', 'http://dbpedia.org/resource/Dimitrios_Kiliaras', 'http://dbpedia.org/resource/Dimitrios_Papadopoulos_(footballer,_born_1950)__Dimitrios_Papadopoulos__1', 'http://dbpedia.org/resource/Dimitrios_Raptakis', 'http://dbpedia.org/resource/Edison_Kola', 'http://dbpedia.org/resource/Eleftherios_Goulielmakis', 'http://dbpedia.org/resource/Emmanuel_Skordilis', 'http://dbpedia.org/resource/Fotis_Kafatos', 'http://dbpedia.org/resource/Franghias_Kavertzas', 'http://dbpedia.org/resource/Georges_Panayotis', 'http://dbpedia.org/resource/Georgios_Kalaitzakis', 'http://dbpedia.org/resource/Georgios_Kastrofylakas', 'http://dbpedia.org/resource/Georgios_Klontzas', 'http://dbpedia.org/resource/Georgios_Kortezas', 'http://dbpedia.org/resource/Georgios_Manousakis', 'http://dbpedia.org/resource/Georgios_Samaras', 'http://dbpedia.org/resource/Georgios_Sarris', 'http://dbpedia.org/resource/Georgios_Sournakis', 'http://dbpedia.org/resource/Georgios_Vlastos', 'http://dbpedia.org/resource/Gianna_Angelopoulos-Daskalaki', 'http://dbpedia.org/resource/Giannis_Boutsakis', 'http://dbpedia.org/resource/Giannis_Dermitzakis', 'http://dbpedia.org/resource/Giannis_Domatas', 'http://dbpedia.org/resource/Giannis_Iatroudis', 'http://dbpedia.org/resource/Giannis_Koutantos', 'http://dbpedia.org/resource/Giannis_Vardinogiannis', 'http://dbpedia.org/resource/Giorgos_Giakoumakis__Giorgos_Giakoumakis__1', 'http://dbpedia.org/resource/Giorgos_Grammatikakis', 'http://dbpedia.org/resource/Giorgos_Petrakis', 'http://dbpedia.org/resource/Ieremias_Palladas', 'http://dbpedia.org/resource/Ilias_Koutsoupias', 'http://dbpedia.org/resource/Ioannis_Apakas', 'http://dbpedia.org/resource/Ioannis_Ikonomou', 'http://dbpedia.org/resource/Ioannis_Pagomenos', 'http://dbpedia.org/resource/Iosif_Daskalakis', 'http://dbpedia.org/resource/Irini_Vasiliou', 'http://dbpedia.org/resource/Joseph_Sifakis', 'http://dbpedia.org/resource/Konstantinos_Gouvis', 'http://dbpedia.org/resource/Konstantinos_Provydakis', 'http://dbpedia.org/resource/Konstanty_Korniakt', 'http://dbpedia.org/resource/Kostas_Chaniotakis', 'http://dbpedia.org/resource/Lampros_Kefaloukos', 'http://dbpedia.org/resource/Lefteris_Gialousis', 'http://dbpedia.org/resource/Lena_Platonos', 'http://dbpedia.org/resource/Leonidas_Kyrkos', 'http://dbpedia.org/resource/Lili_Zografou', 'http://dbpedia.org/resource/Manolis_Bolakis', 'http://dbpedia.org/resource/Manolis_Hatzidakis', 'http://dbpedia.org/resource/Manolis_Moniakis', 'http://dbpedia.org/resource/Manolis_Papamattheakis', 'http://dbpedia.org/resource/Manolis_Patemtzis', 'http://dbpedia.org/resource/Manolis_Roubakis', 'http://dbpedia.org/resource/Manolis_Roussakis', 'http://dbpedia.org/resource/Manolis_Rovithis', 'http://dbpedia.org/resource/Manolis_Saliakas', 'http://dbpedia.org/resource/Manolis_Stefanoudakis', 'http://dbpedia.org/resource/Manuel_Fokas', 'http://dbpedia.org/resource/Markos_Bathas', 'http://dbpedia.org/resource/Meletius_I_Pegas', 'http://dbpedia.org/resource/Michael_Fokas', 'http://dbpedia.org/resource/Michael_Katehakis', 'http://dbpedia.org/resource/Michail_Fragoulakis', 'http://dbpedia.org/resource/Michalis_Sifakis', 'http://dbpedia.org/resource/Michalis_Tzirakis', 'http://dbpedia.org/resource/Minas_Chalkiadakis', 'http://dbpedia.org/resource/Minas_Pitsos', 'http://dbpedia.org/resource/Myron_Michailidis', 'http://dbpedia.org/resource/Nadia_Valavani', 'http://dbpedia.org/resource/Nektarios_Azizi', 'http://dbpedia.org/resource/Nektarios_Tavernarakis', 'http://dbpedia.org/resource/Nicholas_Kalliakis', 'http://dbpedia.org/resource/Nikolaos_Gripiotis', 'http://dbpedia.org/resource/Nikolaos_Ritzos', 'http://dbpedia.org/resource/Nikonas_Spiliotakis', 'http://dbpedia.org/resource/Nikos_Androulakis', 'http://dbpedia.org/resource/Nikos_Karelis', 'http://dbpedia.org/resource/Nikos_Kazantzakis', 'http://dbpedia.org/resource/Nikos_Kousidis', 'http://dbpedia.org/resource/Nikos_Machlas', 'http://dbpedia.org/resource/Nikos_Marinakis', 'http://dbpedia.org/resource/Nikos_Psimopoulos', 'http://dbpedia.org/resource/Nikos_Vrettos', 'http://dbpedia.org/resource/Ntinos_Konstantakis', 'http://dbpedia.org/resource/Odysseas_Elytis', 'http://dbpedia.org/resource/Oresti_Kacurri', 'http://dbpedia.org/resource/Panagiotis_Kalaitzakis', 'http://dbpedia.org/resource/Panagiotis_Petras', 'http://dbpedia.org/resource/Petros_Giakoumakis', 'http://dbpedia.org/resource/Petros_Lambardos', 'http://dbpedia.org/resource/Petros_Marinakis', 'http://dbpedia.org/resource/Şadi_Çalık', 'http://dbpedia.org/resource/Sapfo_Notara', 'http://dbpedia.org/resource/Spyros_Danellis', 'http://dbpedia.org/resource/Stavros_Arnaoutakis', 'http://dbpedia.org/resource/Stavros_Labrakis', 'http://dbpedia.org/resource/Theocharis_Silvestros', 'http://dbpedia.org/resource/Theophanes_the_Cretan', 'http://dbpedia.org/resource/Thomas_Bathas', 'http://dbpedia.org/resource/Vangelis_Nikokyrakis', 'http://dbpedia.org/resource/Vasilios_Babis', 'http://dbpedia.org/resource/Vasilios_Vogiatzis', 'http://dbpedia.org/resource/Victor_(iconographer)', 'http://dbpedia.org/resource/Zacharias_Kavousakis']"
$ws.Range("C8").Value = "['http://dbpedia.org/resource/Anastasios_Triantafyllou_(footballer)', 'http://dbpedia.org/resource/Anastasios_Triantafyllou', 'http://dbpedia.org/resource/Andreas_Bouchalakis', 'http://dbpedia.org/resource/Andreas_Musalus', 'http://dbpedia.org/resource/Antonis_Alexakis', 'http://dbpedia.org/resource/Antonis_Bourselis', 'http://dbpedia.org/resource/Archbishop_Makarios_of_Australia', 'http://dbpedia.org/resource/Beata_Kitsikis', 'http://dbpedia.org/resource/Charalampos_Charalampakis', 'http://dbpedia.org/resource/Charalampos_Damianakis', 'http://dbpedia.org/resource/Christoforos_Liontakis', 'http://dbpedia.org/resource/Christoforos_Papakaliatis', 'http://dbpedia.org/resource/Christos_Belevonis', 'http://dbpedia.org/resource/Christos_Kozyrakis', 'http://dbpedia.org/resource/Christos_Saloustros', 'http://dbpedia.org/resource/Chrysovalantis_Kozoronis', 'http://dbpedia.org/resource/Cyril_Lucaris', 'http://dbpedia.org/resource/Dimitrios_Kiliaras', 'http://dbpedia.org/resource/Dimitrios_Papadopoulos_(footballer,_born_1950)__Dimitrios_Papadopoulos__1', 'http://dbpedia.org/resource/Dimitrios_Raptakis', 'http://dbpedia.org/resource/Edison_Kola', 'http://dbpedia.org/resource/Eleftherios_Goulielmakis', 'http://dbpedia.org/resource/Emmanuel_Skordilis', 'http://dbpedia.org/resource/Fotis_Kafatos', 'http://dbpedia.org/resource/Franghias_Kavertzas', 'http://dbpedia.org/resource/Georges_Panayotis', 'http://dbpedia.org/resource/Georgios_Kalaitzakis', 'http://dbpedia.org/resource/Georgios_Kastrofylakas', 'http://dbpedia.org/resource/Georgios_Klontzas', 'http://dbpedia.org/resource/Georgios_Kortezas', 'http://dbpedia.org/resource/Georgios_Manousakis', 'http://dbpedia.org/resource/Georgios_Samaras', 'http://dbpedia.org/resource/Georgios_Sarris', 'http://dbpedia.org/resource/Georgios_Sournakis', 'http://dbpedia.org/resource/Georgios_Vlastos', 'http://dbpedia.org/resource/Gianna_Angelopoulos-Daskalaki', 'http://dbpedia.org/resource/Giannis_Boutsakis', 'http://dbpedia.org/resource/Giannis_Dermitzakis', 'http://dbpedia.org/resource/Giannis_Domatas', 'http://dbpedia.org/resource/Giannis_Iatroudis', 'http://dbpedia.org/resource/Giannis_Koutantos', 'http://dbpedia.org/resource/Giannis_Vardinogiannis', 'http://dbpedia.org/resource/Giorgos_Giakoumakis__Giorgos_Giakoumakis__1', 'http://dbpedia.org/resource/Giorgos_Grammatikakis', 'http://dbpedia.org/resource/Giorgos_Petrakis', 'http://dbpedia.org/resource/Ieremias_Palladas', 'http://dbpedia.org/resource/Ilias_Koutsoupias', 'http://dbpedia.org/resource/Ioannis_Apakas', 'http://dbpedia.org/resource/Ioannis_Ikonomou', 'http://dbpedia.org/resource/Ioannis_Pagomenos', 'http://dbpedia.org/resource/Iosif_Daskalakis', 'http://dbpedia.org/resource/Irini_Vasiliou', 'http://dbpedia.org/resource/Joseph_Sifakis', 'http://dbpedia.org/resource/Konstantinos_Gouvis', 'http://dbpedia.org/resource/Konstantinos_Provydakis', 'http://dbpedia.org/resource/Konstanty_Korniakt', 'http://dbpedia.org/resource/Kostas_Chaniotakis', 'http://dbpedia.org/resource/Lampros_Kefaloukos', 'http://dbpedia.org/resource/Lefteris_Gialousis', 'http://dbpedia.org/resource/Lena_Platonos', 'http://dbpedia.org/resource/Leonidas_Kyrkos', 'http://dbpedia.org/resource/Lili_Zografou', 'http://dbpedia.org/resource/Manolis_Bolakis', 'http://dbpedia.org/resource/Manolis_Hatzidakis', 'http://dbpedia.org/resource/Manolis_Moniakis', 'http://dbpedia.org/resource/Manolis_Papamattheakis', 'http://dbpedia.org/resource/Manolis_Patemtzis', 'http://dbpedia.org/resource/Manolis_Roubakis', 'http://dbpedia.org/resource/Manolis_Roussakis', 'http://dbpedia.org/resource/Manolis_Rovithis', 'http://dbpedia.org/resource/Manolis_Saliakas', 'http://dbpedia.org/resource/Manolis_Stefanoudakis', 'http://dbpedia.org/resource/Manuel_Fokas', 'http://dbpedia.org/resource/Markos_Bathas', 'http://dbpedia.org/resource/Meletius_I_Pegas', 'http://dbpedia.org/resource/Michael_Fokas', 'http://dbpedia.org/resource/Michael_Katehakis', 'http://dbpedia.org/resource/Michail_Fragoulakis', 'http://dbpedia.org/resource/Michalis_Sifakis', 'http://dbpedia.org/resource/Michalis_Tzirakis', 'http://dbpedia.org/resource/Minas_Chalkiadakis', 'http://dbpedia.org/resource/Minas_Pitsos', 'http://dbpedia.org/resource/Myron_Michailidis', 'http://dbpedia.org/resource/Nadia_Valavani', 'http://dbpedia.org/resource/Nektarios_Azizi', 'http://dbpedia.org/resource/Nektarios_Tavernarakis', 'http://dbpedia.org/resource/Nicholas_Kalliakis', 'http://dbpedia.org/resource/Nikolaos_Gripiotis', 'http://dbpedia.org/resource/Nikolaos_Ritzos', 'http://dbpedia.org/resource/Nikonas_Spiliotakis', 'http://dbpedia.org/resource/Nikos_Androulakis', 'http://dbpedia.org/resource/Nikos_Karelis', 'http://dbpedia.org/resource/Nikos_Kazantzakis', 'http://dbpedia.org/resource/Nikos_Kousidis', 'http://dbpedia.org/resource/Nikos_Machlas', 'http://dbpedia.org/resource/Nikos_Marinakis', 'http://dbpedia.org/resource/Nikos_Psimopoulos', 'http://dbpedia.org/resource/Nikos_Vrettos', 'http://dbpedia.org/resource/Ntinos_Konstantakis', 'http://dbpedia.org/resource/Odysseas_Elytis', 'http://dbpedia.org/resource/Oresti_Kacurri', 'http://dbpedia.org/resource/Panagiotis_Kalaitzakis', 'http://dbpedia.org/resource/Panagiotis_Petras', 'http://dbpedia.org/resource/Petros_Giakoumakis', 'http://dbpedia.org/resource/Petros_Lambardos', 'http://dbpedia.org/resource/Petros_Marinakis', 'http://dbpedia.org/resource/Şadi_Çalık', 'http://dbpedia.org/resource/Sapfo_Notara', 'http://dbpedia.org/resource/Spyros_Danellis', 'http://dbpedia.org/resource/Stavros_Arnaoutakis', 'http://dbpedia.org/resource/Stavros_Labrakis', 'http://dbpedia.org/resource/Theocharis_Silvestros', 'http://dbpedia.org/resource/Theophanes_the_Cretan', 'http://dbpedia.org/resource/Thomas_Bathas', 'http://dbpedia.org/resource/Vangelis_Nikokyrakis', 'http://dbpedia.org/resource/Vasilios_Babis', 'http://dbpedia.org/resource/Vasilios_Vogiatzis', 'http://dbpedia.org/resource/Victor_(iconographer)', 'http://dbpedia.org/resource/Zacharias_Kavousakis']"
$ws.Range("D8").Value = $true

$ws.Range("A9").Value = "What kind of music did Lou Reed play?"
$ws.Range("B9").Value = "['http://dbpedia.org/resource/Art_rock', 'http://dbpedia.org/resource/Experimental_music', 'http://dbpedia.org/resource/Glam_rock', 'http://dbpedia.org/resource/Noise_rock', 'http://dbpedia.org/resource/Proto-punk', 'http://dbpedia.org/resource/Rock_music']"
$ws.Range("C9").Value = "['http://dbpedia.org/resource/Art_rock', 'http://dbpedia.org/resource/Experimental_music', 'http://dbpedia.org/resource/Glam_rock', 'http://dbpedia.org/resource/Noise_rock', 'http://dbpedia.org/resource/Proto-punk', 'http://dbpedia.org/resource/Rock_music']"
$ws.Range("D9").Value = $true

$ws.Range("A10").Value = "Give me the birthdays of all actors of the television show Charmed."
$ws.Range("B10").Value = "['1965-10-01', '1968-07-27', '1969-02-01', '1971-01-26', '1971-04-12', '1972-12-19', '1973-06-15', '1973-09-05', '1973-12-03', '1980-05-19', '1985-11-30']"
$ws.Range("C10").Value = "['1965-10-01', '1968-07-27', '1969-02-01', '1971-01-26', '1971-04-12', '1972-12-19', '1973-06-15', '1973-09-05', '1973-12-03', '1980-05-19', '1985-11-30']"
$ws.Range("D10").Value = $true

$ws.Range("A11").Value = "Give me all people that were born in Vienna and died in Berlin."
$ws.Range("B11").Value = "['http://dbpedia.org/resource/Adolf_Klein', 'http://dbpedia.org/resource/Eduard_von_Winterstein', 'http://dbpedia.org/resource/Eleonore_de_Ahna', 'http://dbpedia.org/resource/Emil_Justitz', 'http://dbpedia.org/resource/Erik_Jan_Hanussen', 'http://dbpedia.org/resource/Eva_Schmidt-Kolmer', 'http://dbpedia.org/resource/Gerda_Laski', 'http://dbpedia.org/resource/Guido_von_List', 'http://dbpedia.org/resource/Hans_Fidesser', 'http://dbpedia.org/resource/Heinrich_Belohlavek', 'http://dbpedia.org/resource/Hugo_Flink', 'http://dbpedia.org/resource/Josephine_Pagay', 'http://dbpedia.org/resource/Karl_Etlinger', 'http://dbpedia.org/resource/Karl_Hellmer', 'http://dbpedia.org/resource/La_Jana_(actress)', 'http://dbpedia.org/resource/Luise_Begas-Parmentier', 'http://dbpedia.org/resource/Mathilde_Danegger', 'http://dbpedia.org/resource/Michael_von_Newlinsky', 'http://dbpedia.org/resource/Richard_Thurnwald', 'http://dbpedia.org/resource/Robert_Bodanzky', 'http://dbpedia.org/resource/Steffie_Spira', 'http://dbpedia.org/resource/Therese_Grünbaum', 'http://dbpedia.org/resource/Tom_Seidmann-Freud']"
$ws.Range("C11").Value = "['http://dbpedia.org/resource/Adolf_Klein', 'http://dbpedia.org/resource/Eduard_von_Winterstein', 'http://dbpedia.org/resource/Eleonore_de_Ahna', 'http://dbpedia.org/resource/Emil_Justitz', 'http://dbpedia.org/resource/Erik_Jan_Hanussen', 'http://dbpedia.org/resource/Eva_Schmidt-Kolmer', 'http://dbpedia.org/resource/Gerda_Laski', 'http://dbpedia.org/resource/Guido_von_List', 'http://dbpedia.org/resource/Hans_Fidesser', 'http://dbpedia.org/resource/Heinrich_Belohlavek', 'http://dbpedia.org/resource/Hugo_Flink', 'http://dbpedia.org/resource/Josephine_Pagay', 'http://dbpedia.org/resource/Karl_Etlinger', 'http://dbpedia.org/resource/Karl_Hellmer', 'http://dbpedia.org/resource/La_Jana_(actress)', 'http://dbpedia.org/resource/Luise_Begas-Parmentier', 'http://dbpedia.org/resource/Mathilde_Danegger', 'http://dbpedia.org/resource/Michael_von_Newlinsky', 'http://dbpedia.org/resource/Richard_Thurnwald', 'http://dbpedia.org/resource/Robert_Bodanzky', 'http://dbpedia.org/resource/Steffie_Spira', 'http://dbpedia.org/resource/Therese_Grünbaum', 'http://dbpedia.org/resource/Tom_Seidmann-Freud']"
$ws.Range("D11").Value = $true

$ws.Range("A12").Value = "Give me a list of all bandleaders that play trumpet."
$ws.Range("B12").Value = "['http://dbpedia.org/resource/Abatte_Barihun', 'http://dbpedia.org/resource/Alan_Silva', 'http://dbpedia.org/resource/Artie_Fields', 'http://dbpedia.org/resource/Charlie_Spivak', 'http://dbpedia.org/resource/Danny_Davis_(country_musician)', 'http://dbpedia.org/resource/Doc_Cheatham', 'http://dbpedia.org/resource/Ernest_%22Doc%22_Paulin', 'http://dbpedia.org/resource/Jack_Jackson_(radio_personality)', 'http://dbpedia.org/resource/Jacobo_Rubalcaba', 'http://dbpedia.org/resource/Jimmy_Dorsey', 'http://dbpedia.org/resource/Kaj_Backlund', 'http://dbpedia.org/resource/Lloyd_Hunter', 'http://dbpedia.org/resource/Marlon_Jordan', 'http://dbpedia.org/resource/Mick_Mulligan', 'http://dbpedia.org/resource/Ollie_Mitchell', 'http://dbpedia.org/resource/Ovie_Alston', 'http://dbpedia.org/resource/Pacho_Galan', 'http://dbpedia.org/resource/Red_Perkins', 'http://dbpedia.org/resource/Richard_Maltby_Sr.', 'http://dbpedia.org/resource/Steamboat_Willie_(musician)']"
$ws.Range("C12").Value = "['http://dbpedia.org/resource/Abatte_Barihun', 'http://dbpedia.org/resource/Alan_Silva', 'http://dbpedia.org/resource/Artie_Fields', 'http://dbpedia.org/resource/Charlie_Spivak', 'http://dbpedia.org/resource/Danny_Davis_(country_musician)', 'http://dbpedia.org/resource/Doc_Cheatham', 'http://dbpedia.org/resource/Ernest_%22Doc%22_Paulin', 'http://dbpedia.org/resource/Jack_Jackson_(radio_personality)', 'http://dbpedia.org/resource/Jacobo_Rubalcaba', 'http://dbpedia.org/resource/Jimmy_Dorsey', 'http://dbpedia.org/resource/Kaj_Backlund', 'http://dbpedia.org/resource/Lloyd_Hunter', 'http://dbpedia.org/resource/Marlon_Jordan', 'http://dbpedia.org/resource/Mick_Mulligan', 'http://dbpedia.org/resource/Ollie_Mitchell', 'http://dbpedia.org/resource/Ovie_Alston', 'http://dbpedia.org/resource/Pacho_Galan', 'http://dbpedia.org/resource/Red_Perkins', 'http://dbpedia.org/resource/Richard_Maltby_Sr.', 'http://dbpedia.org/resource/Steamboat_Willie_(musician)']"
$ws.Range("D12").Value = $true

$ws.Range("A13").Value = "how big is the total area of North Rhine- Westphalia?"
$ws.Range("B13").Value = "['34084130000']"
$ws.Range("C13").Value = "['34084130000']"
$ws.Range("D13").Value = $true
